function Set-TextValue($ws, $ref, $val) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

Set-TextValue $ws "D2" "26.854.72"
Set-TextValue $ws "E2" "  -1.28%  "
Set-TextValue $ws "D3" "1.875.05"
Set-TextValue $ws "E3" "  -1.60%  "
Set-TextValue $ws "D5" "301.44"
Set-TextValue $ws "E5" "  -2.05%  "
Set-TextValue $ws "E6" "  -0.15%  "
Set-TextValue $ws "D7" "0.5369"
Set-TextValue $ws "E7" "  +1.84%  "
Set-TextValue $ws "D8" "0.3753"
Set-TextValue $ws "E8" "  -1.91%  "
Set-TextValue $ws "D9" "0.07185"
Set-TextValue $ws "E9" "  -1.62%  "
Set-TextValue $ws "D10" "21.61"
Set-TextValue $ws "E10" "  +0.27%  "
Set-TextValue $ws "D11" "0.8898"
Set-TextValue $ws "E11" "  -1.69%  "
Set-TextValue $ws "D12" "0.08150"
Set-TextValue $ws "E12" "  +0.62%  "
Set-TextValue $ws "B13" "WrappedEther"
Set-TextValue $ws "C13" "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue $ws "D13" "1.887.25"
Set-TextValue $ws "E13" "  +2.79%  "
Set-TextValue $ws "B14" "Litecoin"
Set-TextValue $ws "C14" "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextValue $ws "D14" "93.59"
Set-TextValue $ws "E14" "  -2.47%  "
Set-TextValue $ws "D15" "5.312"
Set-TextValue $ws "E15" "  -1.15%  "
Set-TextValue $ws "E16" "  -0.04%  "
Set-TextValue $ws "D17" "14.84"
Set-TextValue $ws "E17" "  +0.69%  "
Set-TextValue $ws "D18" "0.000008550"
Set-TextValue $ws "E18" "  -1.48%  "
Set-TextValue $ws "E19" "  -0.17%  "
Set-TextValue $ws "D20" "26.889.97"
Set-TextValue $ws "E20" "  -1.28%  "
Set-TextValue $ws "E21" "  -2.63%  "
Set-TextValue $ws "D23" "6.406"
Set-TextValue $ws "E23" "  -1.36%  "
Set-TextValue $ws "D24" "2.304"
Set-TextValue $ws "E24" "  -1.62%  "
Set-TextValue $ws "D25" "146.35"
Set-TextValue $ws "E25" "  -2.58%  "
Set-TextValue $ws "D26" "18.06"
Set-TextValue $ws "E26" "  -1.07%  "
Set-TextValue $ws "D27" "1.731"
Set-TextValue $ws "E27" "  -0.72%  "
Set-TextValue $ws "D28" "113.94"
Set-TextValue $ws "E28" "  -2.63%  "
Set-TextValue $ws "D29" "4.730"
Set-TextValue $ws "E29" "  -2.37%  "
Set-TextValue $ws "D30" "4.620"
Set-TextValue $ws "E30" "  -5.15%  "
Set-TextValue $ws "D31" "0.09158"
Set-TextValue $ws "E31" "  -0.76%  "
Set-TextValue $ws "D32" "0.8176"
Set-TextValue $ws "E32" "  +0.68%  "
Set-TextValue $ws "D33" "0.05003"
Set-TextValue $ws "E33" "  -1.27%  "
Set-TextValue $ws "D34" "1.176"
Set-TextValue $ws "E34" "  -4.43%  "
Set-TextValue $ws "D35" "2.949"
Set-TextValue $ws "E35" "  -1.20%  "
Set-TextValue $ws "D36" "0.6063"
Set-TextValue $ws "E36" "  +5.75%  "
Set-TextValue $ws "D37" "3.220"
Set-TextValue $ws "E37" "  -4.08%  "
Set-TextValue $ws "E38" "  -3.03%  "
Set-TextValue $ws "D39" "0.01953"
Set-TextValue $ws "E39" "  -2.09%  "
Set-TextValue $ws "D40" "1.070"
Set-TextValue $ws "E40" "  -1.59%  "
Set-TextValue $ws "D41" "6.640"
Set-TextValue $ws "E41" "  +0.29%  "
Set-TextValue $ws "E42" "  -0.71%  "
Set-TextValue $ws "D43" "115.10"
Set-TextValue $ws "E43" "  -1.48%  "
Set-TextValue $ws "D44" "0.5111"
Set-TextValue $ws "E44" "  +3.67%  "
Set-TextValue $ws "E45" "  -1.70%  "
Set-TextValue $ws "E46" "  -0.14%  "
Set-TextValue $ws "D47" "9.950"
Set-TextValue $ws "E47" "  -2.23%  "
Set-TextValue $ws "E48" "  -0.30%  "
Set-TextValue $ws "D49" "37.77"
Set-TextValue $ws "E49" "  -2.00%  "
Set-TextValue $ws "D50" "0.06085"
Set-TextValue $ws "E50" "  +2.02%  "
Set-TextValue $ws "D51" "62.23"
Set-TextValue $ws "E51" "  -3.25%  "
